# Rename the inline logo pictures in the document's headers and footers.
#   Headers: BTec_Logo-Orange  image1.jpg -> image2.jpg
#   Footers: PearsonLogo.png   image2.png -> image1.png
$d = $word.ActiveDocument
$sec = $d.Sections(1)

for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.Name -eq "image1.jpg") {
                $shp.Name = "image2.jpg"
            }
        }
    }

    $ftr = $sec.Footers($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.Name -eq "image2.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
